$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04794271460666266
$ws.Range("D2").Value = 0.193508410118568
$ws.Range("E2").Value = 0.06606112079597182
$ws.Range("F2").Value = 3.54499807292251
$ws.Range("G2").Value = 0.002509703061337871
$ws.Range("K2").Value = 1.517009757728829
$ws.Range("M2").Value = 0.4448321100705996
$ws.Range("C3").Value = 0.04260015905781245
$ws.Range("D3").Value = 0.1833137481967668
$ws.Range("E3").Value = 0.06027124878698586
$ws.Range("F3").Value = 3.339990571879412
$ws.Range("G3").Value = 0.002516736769535932
$ws.Range("K3").Value = 1.403504704817351
$ws.Range("M3").Value = 0.4099697860196798
$ws.Range("C4").Value = 0.0393426102755825
$ws.Range("D4").Value = 0.177068457234455
$ws.Range("E4").Value = 0.05674559003634627
$ws.Range("F4").Value = 3.215251527776445
$ws.Range("G4").Value = 0.002521271152526035
$ws.Range("K4").Value = 1.33531791114072
$ws.Range("M4").Value = 0.3889290001933148
$ws.Range("C5").Value = 0.03802064977891462
$ws.Range("D5").Value = 0.1745263476503283
$ws.Range("E5").Value = 0.05531593101532195
$ws.Range("F5").Value = 3.164694216524197
$ws.Range("G5").Value = 0.002523173408698779
$ws.Range("K5").Value = 1.307904775614162
$ws.Range("M5").Value = 0.3804448538546978
$ws.Range("C6").Value = 0.03780146620368896
$ws.Range("D6").Value = 0.1741043871542303
$ws.Range("E6").Value = 0.05507895576092636
$ws.Range("F6").Value = 3.15631549847447
$ws.Range("G6").Value = 0.002523492572746546
$ws.Range("K6").Value = 1.303375259357892
$ws.Range("M6").Value = 0.379041466415174
$ws.Range("C7").Value = 0.03932475979848959
$ws.Range("D7").Value = 0.1770341626363461
$ws.Range("E7").Value = 0.05672628090525933
$ws.Range("F7").Value = 3.214568595598024
$ws.Range("G7").Value = 0.002521296586350363
$ws.Range("K7").Value = 1.334946702216712
$ws.Range("M7").Value = 0.3888142171165327
$ws.Range("C8").Value = 0.04609573591413607
$ws.Range("D8").Value = 0.1899899079802339
$ws.Range("E8").Value = 0.0640584966585962
$ws.Range("F8").Value = 3.474068977664729
$ws.Range("G8").Value = 0.002512083662825871
$ws.Range("K8").Value = 1.477557547601634
$ws.Range("M8").Value = 0.4327348824909123
$ws.Range("C9").Value = 0.05956507685074541
$ws.Range("D9").Value = 0.2155416957839122
$ws.Range("E9").Value = 0.07868429503253793
$ws.Range("F9").Value = 3.992484918094618
$ws.Range("G9").Value = 0.002495717783382024
$ws.Range("K9").Value = 1.769416227651391
$ws.Range("M9").Value = 0.5218357450691684
$ws.Range("C10").Value = 0.06959348479192329
$ws.Range("D10").Value = 0.2344498723364268
$ws.Range("E10").Value = 0.08960189360917781
$ws.Range("F10").Value = 4.379939922884034
$ws.Range("G10").Value = 0.002484715731371429
$ws.Range("K10").Value = 1.991661753841868
$ws.Range("M10").Value = 0.5892269358565301
$ws.Range("C11").Value = 0.07418797172843483
$ws.Range("D11").Value = 0.2430913140926236
$ws.Range("E11").Value = 0.09461066809488017
$ws.Range("F11").Value = 4.557800156230428
$ws.Range("G11").Value = 0.002479929333560079
$ws.Range("K11").Value = 2.094550732736423
$ws.Range("M11").Value = 0.6203300863466268
$ws.Range("C12").Value = 0.07593273458562066
$ws.Range("D12").Value = 0.2463701724324778
$ws.Range("E12").Value = 0.09651381478889931
$ws.Range("F12").Value = 4.625395293311612
$ws.Range("G12").Value = 0.002478148020288323
$ws.Range("K12").Value = 2.133776058968238
$ws.Range("M12").Value = 0.6321743186085627
$ws.Range("C13").Value = 0.07555674528603618
$ws.Range("D13").Value = 0.2456637086179398
$ws.Range("E13").Value = 0.09610364638171376
$ws.Range("F13").Value = 4.610826447114732
$ws.Range("G13").Value = 0.002478530274391277
$ws.Range("K13").Value = 2.125316360824172
$ws.Range("M13").Value = 0.6296204782224635
$ws.Range("C14").Value = 0.07433141412501243
$ws.Range("D14").Value = 0.2433609323365715
$ws.Range("E14").Value = 0.09476711017808981
$ws.Range("F14").Value = 4.563356299618533
$ws.Range("G14").Value = 0.002479782160092291
$ws.Range("K14").Value = 2.097772505962496
$ws.Range("M14").Value = 0.6213031815995862
$ws.Range("C15").Value = 0.07358151315382599
$ws.Range("D15").Value = 0.2419512892907676
$ws.Range("E15").Value = 0.093949291338987
$ws.Range("F15").Value = 4.534311527721741
$ws.Range("G15").Value = 0.002480553031424636
$ws.Range("K15").Value = 2.080935606740013
$ws.Range("M15").Value = 0.6162172693285015
$ws.Range("C16").Value = 0.06929390018628112
$ws.Range("D16").Value = 0.2338860052469443
$ws.Range("E16").Value = 0.0892754393910522
$ws.Range("F16").Value = 4.368349700961573
$ws.Range("G16").Value = 0.00248503291086338
$ws.Range("K16").Value = 1.984974241123211
$ws.Range("M16").Value = 0.5872034215242081
$ws.Range("C17").Value = 0.06667209563741494
$ws.Range("D17").Value = 0.2289490137497694
$ws.Range("E17").Value = 0.08641926417762136
$ws.Range("F17").Value = 4.266957419007724
$ws.Range("G17").Value = 0.002487836971054906
$ws.Range("K17").Value = 1.926567545999092
$ws.Range("M17").Value = 0.5695200971510701
$ws.Range("C18").Value = 0.06516714776611821
$ws.Range("D18").Value = 0.2261131062799109
$ws.Range("E18").Value = 0.08478042833534971
$ws.Range("F18").Value = 4.208789568190753
$ws.Range("G18").Value = 0.00248947037135564
$ws.Range("K18").Value = 1.893141499546289
$ws.Range("M18").Value = 0.559391030537526
$ws.Range("C19").Value = 0.06465811354667039
$ws.Range("D19").Value = 0.2251535320619951
$ws.Range("E19").Value = 0.08422621638071348
$ws.Range("F19").Value = 4.189120370911098
$ws.Range("G19").Value = 0.002490026954100814
$ws.Range("K19").Value = 1.881852668380873
$ws.Range("M19").Value = 0.5559686434110773
$ws.Range("C20").Value = 0.06695087425708834
$ws.Range("D20").Value = 0.229474174376378
$ws.Range("E20").Value = 0.08672289648875875
$ws.Range("F20").Value = 4.277735154178572
$ws.Range("G20").Value = 0.002487536345973029
$ws.Range("K20").Value = 1.932767611977454
$ws.Range("M20").Value = 0.5713981654786835
$ws.Range("C21").Value = 0.07469118772013417
$ws.Range("D21").Value = 0.2440371300255038
$ws.Range("E21").Value = 0.09515950586816047
$ws.Range("F21").Value = 4.577292727815632
$ws.Range("G21").Value = 0.002479413606278121
$ws.Range("K21").Value = 2.10585560096365
$ws.Range("M21").Value = 0.6237443624669226
$ws.Range("C22").Value = 0.07977879050326919
$ws.Range("D22").Value = 0.2535933053540873
$ws.Range("E22").Value = 0.1007109780699054
$ws.Range("F22").Value = 4.77449519497776
$ws.Range("G22").Value = 0.002474286628953603
$ws.Range("K22").Value = 2.220517582425941
$ws.Range("M22").Value = 0.6583420994964797
$ws.Range("C23").Value = 0.07706071465536013
$ws.Range("D23").Value = 0.2484892174161075
$ws.Range("E23").Value = 0.09774448988449791
$ws.Range("F23").Value = 4.66911000531519
$ws.Range("G23").Value = 0.002477006442614239
$ws.Range("K23").Value = 2.15917736491906
$ws.Range("M23").Value = 0.6398406230615592
$ws.Range("C24").Value = 0.06682483110319026
$ws.Range("D24").Value = 0.2292367417322794
$ws.Range("E24").Value = 0.08658561432593359
$ws.Range("F24").Value = 4.272862155060949
$ws.Range("G24").Value = 0.002487672192154316
$ws.Range("K24").Value = 1.929964086872189
$ws.Range("M24").Value = 0.5705489747141002
$ws.Range("C25").Value = 0.05589915484209484
$ws.Range("D25").Value = 0.2086088988375536
$ws.Range("E25").Value = 0.07469897424857663
$ws.Range("F25").Value = 3.85114077886746
$ws.Range("G25").Value = 0.002499964650560755
$ws.Range("K25").Value = 1.689116059873072
$ws.Range("M25").Value = 0.4974014779592295
